$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) cells we are updating to Text format so that
# numeric-looking strings (e.g. "4.511", "0.000008500") are preserved
# exactly as text instead of being parsed/normalized as numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D21", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Updated Price values
$ws.Range("D2").Value = "26.042.75"
$ws.Range("D3").Value = "1.666.77"
$ws.Range("D5").Value = "215.84"
$ws.Range("D6").Value = "0.5093"
$ws.Range("D8").Value = "0.2664"
$ws.Range("D11").Value = "0.07462"
$ws.Range("D12").Value = "1.684.17"
$ws.Range("D13").Value = "4.511"
$ws.Range("D14").Value = "0.5803"
$ws.Range("D15").Value = "0.000008500"
$ws.Range("D16").Value = "64.02"
$ws.Range("D17").Value = "26.126.86"
$ws.Range("D18").Value = "4.917"
$ws.Range("D21").Value = "189.77"
$ws.Range("D24").Value = "144.99"
$ws.Range("D25").Value = "7.595"
$ws.Range("D26").Value = "0.1205"
$ws.Range("D28").Value = "0.06613"
$ws.Range("D29").Value = "1.329"
$ws.Range("D31").Value = "3.551"
$ws.Range("D32").Value = "3.515"
$ws.Range("D33").Value = "1.657"
$ws.Range("D35").Value = "0.6138"
$ws.Range("D36").Value = "2.368"
$ws.Range("D37").Value = "2.686"
$ws.Range("D38").Value = "6.369"
$ws.Range("D39").Value = "1.091.84"
$ws.Range("D41").Value = "0.8691"
$ws.Range("D43").Value = "101.30"
$ws.Range("D44").Value = "1.813.18"
$ws.Range("D45").Value = "0.00000000115"
$ws.Range("D47").Value = "1.009"
$ws.Range("D48").Value = "8.053"
$ws.Range("D49").Value = "0.05229"
$ws.Range("D51").Value = "5.997"

# Updated Volume(1h) values
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +2.60%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").Value = "  +14.16%  "
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +8.22%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("E51").Value = "  +2.98%  "
